# Update the "Students" header row: rename the "Name" column header to "Full Name".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Students")

$ws.Range("B1").Value = "Full Name"

# Move/restore the active cell selection as recorded in the sheet view.
$ws.Range("E21").Select()
